$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "Daily Achievement" column (N) one column to the left (M),
# carrying the custom column width along with it.
$ws.Columns.Item(14).Cut()
$ws.Columns.Item(13).Insert()

# Re-establish the formula (Cut/Insert already preserves it, but make sure).
$ws.Range("M5").Formula = "=L5-L4"

# New day of writing - row 6
$ws.Range("A6").Value = (Get-Date -Year 2012 -Month 6 -Day 29)
$ws.Range("A6").NumberFormat = "m/d/yyyy"
$ws.Range("B6").Value = 6686
$ws.Range("C6").Value = 11526
$ws.Range("D6").Value = 5327
$ws.Range("E6").Value = 2585
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 10
$ws.Range("H6").Value = 253
$ws.Range("I6").Value = 128
$ws.Range("J6").Value = 7
$ws.Range("K6").Formula = "=SUM(B6:J6)"
$ws.Range("L6").Value = 26968
$ws.Range("M6").Formula = "=L6-L5"
$ws.Range("N6").Value = "Returned to writing after long break with illness -small word count but good writing"

$ws.Range("N12").Select()
